{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: the run that hosts the inline picture gets \"do not spell check\"\n//           turned on (OOXML: <w:noProof/> added to that run's <w:rPr>).\n// Change 2: seventeen new paragraphs listing the AOI \"Color Values\" legend\n//           (\"Color Values:\" + 16 numbered color entries) are inserted right\n//           after the empty paragraph that follows the \"Index 80 ...\"\n//           paragraph, and right before the document's trailing empty\n//           paragraph.\n\nconst body = context.document.body;\n\n// --- Change 1: mark the picture's run as \"do not spell check\" (noProof) ---\nconst inlinePictures = body.inlinePictures;\ninlinePictures.load(\"items\");\nawait context.sync();\n\nfor (const pic of inlinePictures.items) {\n  const pictureRange = pic.getRange();\n  pictureRange.hasNoProofing = true;\n}\nawait context.sync();\n\n// --- Change 2: insert the \"Color Values\" legend paragraphs ---\n// Locate the empty paragraph that immediately follows the\n// \"Index 80 \u2013 Must be = 2 ...\" paragraph; the new content goes right after\n// it (and therefore right before the document's final empty paragraph).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet indexParaPos = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Index 80\") !== -1) {\n    indexParaPos = i;\n    break;\n  }\n}\n\nif (indexParaPos === -1) {\n  throw new Error(\"Could not find the 'Index 80' paragraph.\");\n}\n\nconst anchorPara = paragraphs.items[indexParaPos + 1];\n\nconst colorLines = [\n  \"Color Values: \",\n  \"0 = Green \",\n  \"1 = Red \",\n  \"2 = Orange \",\n  \"3 = Amber \",\n  \"4 = Yellow \",\n  \"5 = Lime Green \",\n  \"6 = Spring Green \",\n  \"7 = Cyan \",\n  \"8 = Sky Blue \",\n  \"9 = Blue \",\n  \"10 = Violet \",\n  \"11 = Magenta \",\n  \"12 = Rose \",\n  \"13 = White \",\n  \"14 = Custom1 \",\n  \"15 = Custom2\",\n];\n\nlet insertAfter = anchorPara;\nfor (const line of colorLines) {\n  const newPara = insertAfter.insertParagraph(line, \"After\");\n  newPara.font.size = 14; // sz/szCs 28 half-points == 14pt, matches surrounding text\n  insertAfter = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: the run that hosts the inline picture gets NoProofing turned on\n#           (OOXML: <w:noProof/> added to that run's <w:rPr>).\n# Change 2: sixteen new paragraphs listing the AOI \"Color Values\" legend are\n#           inserted right after the empty paragraph that follows the\n#           \"Index 80 ...\" paragraph, and right before the document's\n#           trailing empty paragraph.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: mark the picture's range as \"do not spell check\" (NoProofing) ---\nif ($d.InlineShapes.Count -ge 1) {\n    $shape = $d.InlineShapes(1)\n    $shape.Range.NoProofing = 1\n}\n\n# --- Change 2: insert the \"Color Values\" legend paragraphs ---\n# Find the paragraph that contains \"Index 80\"; the new content is inserted\n# right after the (empty) paragraph that follows it.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -match \"Index 80\") {\n        $anchorIndex = $i + 1\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'Index 80' paragraph.\"\n}\n\n$lines = @(\n    \"Color Values: \",\n    \"0 = Green \",\n    \"1 = Red \",\n    \"2 = Orange \",\n    \"3 = Amber \",\n    \"4 = Yellow \",\n    \"5 = Lime Green \",\n    \"6 = Spring Green \",\n    \"7 = Cyan \",\n    \"8 = Sky Blue \",\n    \"9 = Blue \",\n    \"10 = Violet \",\n    \"11 = Magenta \",\n    \"12 = Rose \",\n    \"13 = White \",\n    \"14 = Custom1 \",\n    \"15 = Custom2\"\n)\n\n$insertAfterIndex = $anchorIndex\nforeach ($line in $lines) {\n    $p = $d.Paragraphs($insertAfterIndex)\n    $p.Range.InsertParagraphAfter()\n    $insertAfterIndex = $insertAfterIndex + 1\n    $newPara = $d.Paragraphs($insertAfterIndex)\n    $newPara.Range.Text = $line\n}\n"}
